$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 7 - India
$ws.Range("D7").Value = 154330
$ws.Range("E7").Value = 146383

# Row 18 - Pakistan
$ws.Range("B18").Value = 132405
$ws.Range("C18").Value = 6472
$ws.Range("D18").Value = 50056
$ws.Range("E18").Value = 79798
$ws.Range("G18").Value = 88
$ws.Range("H18").Value = 2551

# Row 56 - Kazajistan
$ws.Range("B56").Value = 14238
$ws.Range("C56").Value = 366
$ws.Range("E56").Value = 5339

# Row 71 - Australia
$ws.Range("B71").Value = 7294
$ws.Range("C71").Value = 4
$ws.Range("D71").Value = 6803
$ws.Range("E71").Value = 389

# Row 99 - Kirguistan
$ws.Range("B99").Value = 2207
$ws.Range("C99").Value = 41
$ws.Range("D99").Value = 1722
$ws.Range("E99").Value = 458
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 27

# Row 183 - Butan
$ws.Range("D183").Value = 20
$ws.Range("E183").Value = 42
